# rill-analysis: Enhance browser compatible
#
# The "trend" sheet shows, for each hour row (28-51), both a percentage
# (ratio) change and an absolute change versus the comparison periods
# defined on the _input sheet. The U and AC columns used to surface the
# *ratio* figures (_input!F / _input!G) but should instead surface the
# *absolute difference* figures (_input!H / _input!I) - matching what
# columns N and V already show - so that the numbers render consistently
# across browsers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trend")

for ($row = 28; $row -le 51; $row++) {
    $inputRow = $row - 23

    $ws.Range("U$row").Formula = "=_input!`$H$inputRow"
    $ws.Range("AC$row").Formula = "=_input!`$I$inputRow"
}
